$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "20.246.52"
$ws.Cells.Item(2, 5).Value = "  +1.53%  "
$ws.Cells.Item(3, 4).Value = "1.440.05"
$ws.Cells.Item(3, 5).Value = "  +2.12%  "
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.010"
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.81%  "
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9147"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -8.64%  "
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "277.63"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.82%  "
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.3664"
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.57%  "
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.3127"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +1.79%  "
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "39.01"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.61%  "
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.019"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +2.59%  "
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06524"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.36%  "
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.12%  "
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.393"
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.37%  "
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.52"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +3.51%  "
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.077"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -1.04%  "
$ws.Cells.Item(16, 4).Value = "1.446.25"
$ws.Cells.Item(16, 5).Value = "  +2.53%  "
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.00001017"
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.00%  "
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9266"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -7.43%  "
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.05619"
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -2.24%  "
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "67.21"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -8.25%  "
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.412"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -3.17%  "
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "14.40"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.90"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.85%  "
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.258"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -1.09%  "
$ws.Cells.Item(25, 4).Value = "20.289.88"
$ws.Cells.Item(25, 5).Value = "  +1.70%  "
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.190"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -2.28%  "
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "135.10"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -2.57%  "
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "16.96"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.23%  "
$ws.Cells.Item(29, 4).Value = "1.597.05"
$ws.Cells.Item(29, 5).Value = "  +1.80%  "
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "110.25"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +1.30%  "
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.692"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -2.73%  "
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.8153"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.27%  "
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.872"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -8.32%  "
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.07636"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.50%  "
$ws.Cells.Item(35, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.486"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +16.55%  "
$ws.Cells.Item(36, 2).Value = "Hedera"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.05967"
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +3.92%  "
$ws.Cells.Item(37, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.680"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.35%  "
$ws.Cells.Item(38, 2).Value = "TrustWalletToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.141"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +7.82%  "
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.24"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -0.20%  "
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.01990"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -2.12%  "
$ws.Cells.Item(41, 2).Value = "Frax"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9283"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -7.23%  "
$ws.Cells.Item(42, 2).Value = "Algorand"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1820"
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -6.05%  "
$ws.Cells.Item(43, 2).Value = "FraxShare"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.961"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -17.45%  "
$ws.Cells.Item(44, 2).Value = "PancakeSwap"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.521"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.07%  "
$ws.Cells.Item(45, 2).Value = "TheSandbox"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.5220"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.06%  "
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.96"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.26%  "
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "120.02"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +8.16%  "
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.5135"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +0.76%  "
$ws.Cells.Item(49, 5).Value = "  -1.93%  "
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06330"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +2.80%  "
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9949"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.59%  "
